# Generate Report for Handback
# Marks the zh-cn and de-de handoffs as handed back: updates the Status
# column, stamps the handback datetime, and records the "Latest Target
# File" / "Latest Handback File" hyperlinks (columns F and G) that were
# previously empty.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn = $wb.Worksheets.Item("zh-cn")
$dede = $wb.Worksheets.Item("de-de")

$statusText = "Handed back: in sync with en-US"
$zhHandbackTime = "2016-03-23 03:06:24"
$deHandbackTime = "2016-03-23 03:06:30"

$mdDisplay = "ef19ff1f-a741-4920-afa9-d4358ee2485d.md"
$zhXlfDisplay = "ef19ff1f-a741-4920-afa9-d4358ee2485d.bebee68a9ca2ee6b78adf27937e488b8b7425634.zh-cn.xlf"
$deXlfDisplay = "ef19ff1f-a741-4920-afa9-d4358ee2485d.bebee68a9ca2ee6b78adf27937e488b8b7425634.de-de.xlf"

# Grab the existing hyperlink targets for the "ef19ff1f..." entry (row 2) on
# each language sheet so the new "Latest Target File" / "Latest Handback
# File" links point at the same place as the existing "Source File Name" /
# "Latest Handoff File" links.
$zhMdUrl = $zhcn.Range("A2").Hyperlinks.Item(1).Address
$zhXlfUrl = $zhcn.Range("D2").Hyperlinks.Item(1).Address
$deMdUrl = $dede.Range("A2").Hyperlinks.Item(1).Address
$deXlfUrl = $dede.Range("D2").Hyperlinks.Item(1).Address

# --- Overview sheet: both language statuses flip to "handed back" too ---
$overview.Range("B2").Value = $statusText
$overview.Range("C2").Value = $statusText
$overview.Range("B3").Value = $statusText
$overview.Range("C3").Value = $statusText

# --- zh-cn sheet ---
foreach ($row in 2,3) {
    $zhcn.Range("C$row").Value = $statusText

    $zhcn.Range("F$row").Value = $mdDisplay
    $zhcn.Hyperlinks.Add($zhcn.Range("F$row"), $zhMdUrl, "", "", $mdDisplay) | Out-Null

    $zhcn.Range("G$row").Value = $zhXlfDisplay
    $zhcn.Hyperlinks.Add($zhcn.Range("G$row"), $zhXlfUrl, "", "", $zhXlfDisplay) | Out-Null

    $zhcn.Range("H$row").Value = $zhHandbackTime
}

# --- de-de sheet ---
foreach ($row in 2,3) {
    $dede.Range("C$row").Value = $statusText

    $dede.Range("F$row").Value = $mdDisplay
    $dede.Hyperlinks.Add($dede.Range("F$row"), $deMdUrl, "", "", $mdDisplay) | Out-Null

    $dede.Range("G$row").Value = $deXlfDisplay
    $dede.Hyperlinks.Add($dede.Range("G$row"), $deXlfUrl, "", "", $deXlfDisplay) | Out-Null

    $dede.Range("H$row").Value = $deHandbackTime
}
